$d = $word.ActiveDocument

$replacements = @(
    @{old="69×40=2760"; new="23×46=1058"},
    @{old="97×60=5820"; new="29×18=522"},
    @{old="72×29=2088"; new="78×71=5538"},
    @{old="77×91=7007"; new="17×26=442"},
    @{old="71×40=2840"; new="40×52=2080"},
    @{old="45×60=2700"; new="69×46=3174"},
    @{old="18×21=378";  new="41×51=2091"},
    @{old="40×95=3800"; new="90×65=5850"},
    @{old="33×93=3069"; new="79×67=5293"},
    @{old="22×39=858";  new="54×50=2700"},
    @{old="55×39=2145"; new="50×18=900"},
    @{old="43×17=731";  new="46×96=4416"},
    @{old="67×69=4623"; new="54×24=1296"},
    @{old="94×85=7990"; new="98×22=2156"},
    @{old="98×54=5292"; new="89×16=1424"},
    @{old="29×61=1769"; new="13×97=1261"},
    @{old="40×75=3000"; new="96×52=4992"},
    @{old="39×63=2457"; new="70×94=6580"},
    @{old="48×71=3408"; new="49×81=3969"},
    @{old="17×63=1071"; new="86×28=2408"},
    @{old="93×24=2232"; new="72×15=1080"},
    @{old="73×87=6351"; new="59×35=2065"},
    @{old="60×45=2700"; new="82×94=7708"},
    @{old="29×55=1595"; new="21×50=1050"},
    @{old="34×76=2584"; new="15×12=180"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
